# add support to inputfield
# Mark the "completion status" column (F) with "V" for the rows that now
# have support: InputField (row 10) in the upper table, and image /
# slicedImage / texture / text (rows 15, 16, 17, 22) in the lower table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F10").Value = "V"
$ws.Range("F15").Value = "V"
$ws.Range("F16").Value = "V"
$ws.Range("F17").Value = "V"
$ws.Range("F22").Value = "V"

# Move the active selection to F10, matching the cell that was just edited.
$ws.Range("F10").Select()
